# Se resta el cupo de la clase anotada
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the "Nombre" data in a new column D, matching each row's email address.
$ws.Range("D1").Value = "PEREZ JUAN"
$ws.Range("D2").Value = "GIMENEZ ZAIRA"
$ws.Range("D3").Value = "SUAREZ JULIETA"
$ws.Range("D4").Value = "SANCHEZ ROBERTO"

# Update the active selection to reflect where the cursor ended up.
$ws.Range("D8").Select()
